$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ltf"
$ws.Range("C2").Value = "Lrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.245485
$ws.Range("H2").Value = 0.736455
$ws.Range("I2").Value = 0.01511172246591349
$ws.Range("J2").Value = 0.01511172246591349
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.749051
$ws.Range("N2").Value = 5.247153000000001
$ws.Range("O2").Value = 0.003644723415756578
$ws.Range("P2").Value = 0.003644723415756579
$ws.Range("Q2").Value = 0.429365784735
$ws.Range("R2").Value = 3.864292062615001
$ws.Range("S2").Value = [double]"5.507804872392965e-05"
$ws.Range("T2").Value = [double]"5.507804872392965e-05"

$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ltf"
$ws.Range("C3").Value = "Lrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.245485
$ws.Range("H3").Value = 0.736455
$ws.Range("I3").Value = 0.01511172246591349
$ws.Range("J3").Value = 0.01511172246591349
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 143.0355533333334
$ws.Range("N3").Value = 429.10666
$ws.Range("O3").Value = 0.2980616520156925
$ws.Range("P3").Value = 0.2980616520156925
$ws.Range("Q3").Value = 35.11308281003333
$ws.Range("R3").Value = 316.0177452903
$ws.Range("S3").Value = 0.00450422496299283
$ws.Range("T3").Value = 0.004504224962992829

$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ltf"
$ws.Range("C4").Value = "Lrp1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.245485
$ws.Range("H4").Value = 0.736455
$ws.Range("I4").Value = 0.01511172246591349
$ws.Range("J4").Value = 0.01511172246591349
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 132.804812
$ws.Range("N4").Value = 398.414436
$ws.Range("O4").Value = 0.2767425352500014
$ws.Range("P4").Value = 0.2767425352500014
$ws.Range("Q4").Value = 32.60158927382
$ws.Range("R4").Value = 293.41430346438
$ws.Range("S4").Value = 0.004182056387211302
$ws.Range("T4").Value = 0.004182056387211303

$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Ltf"
$ws.Range("C5").Value = "Lrp1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.245485
$ws.Range("H5").Value = 0.736455
$ws.Range("I5").Value = 0.01511172246591349
$ws.Range("J5").Value = 0.01511172246591349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 125.707184
$ws.Range("N5").Value = 377.121552
$ws.Range("O5").Value = 0.261952291301752
$ws.Range("P5").Value = 0.261952291301752
$ws.Range("Q5").Value = 30.85922806424
$ws.Range("R5").Value = 277.73305257816
$ws.Range("S5").Value = 0.003958550325462201
$ws.Range("T5").Value = 0.003958550325462201

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ltf"
$ws.Range("C6").Value = "Lrp1"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.245485
$ws.Range("H6").Value = 0.736455
$ws.Range("I6").Value = 0.01511172246591349
$ws.Range("J6").Value = 0.01511172246591349
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 53.58648166666666
$ws.Range("N6").Value = 160.759445
$ws.Range("O6").Value = 0.111665071229204
$ws.Range("P6").Value = 0.1116650712292041
$ws.Range("Q6").Value = 13.15467745194167
$ws.Range("R6").Value = 118.392097067475
$ws.Range("S6").Value = 0.001687451565552193
$ws.Range("T6").Value = 0.001687451565552193

$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ltf"
$ws.Range("C7").Value = "Lrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.245485
$ws.Range("H7").Value = 0.736455
$ws.Range("I7").Value = 0.01511172246591349
$ws.Range("J7").Value = 0.01511172246591349
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.00271466666667
$ws.Range("N7").Value = 69.008144
$ws.Range("O7").Value = 0.04793372678759353
$ws.Range("P7").Value = 0.04793372678759355
$ws.Range("Q7").Value = 5.646821409946666
$ws.Range("R7").Value = 50.82139268952
$ws.Range("S7").Value = 0.0007243611759710367
$ws.Range("T7").Value = 0.0007243611759710368

$ws.Range("A8").Value = "Neutro"
$ws.Range("B8").Value = "Ltf"
$ws.Range("C8").Value = "Lrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.99918866666667
$ws.Range("H8").Value = 47.99756600000001
$ws.Range("I8").Value = 0.9848882775340866
$ws.Range("J8").Value = 0.9848882775340865
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.749051
$ws.Range("N8").Value = 5.247153000000001
$ws.Range("O8").Value = 0.003644723415756578
$ws.Range("P8").Value = 0.003644723415756579
$ws.Range("Q8").Value = 27.98339693662201
$ws.Range("R8").Value = 251.8505724295981
$ws.Range("S8").Value = 0.003589645367032649
$ws.Range("T8").Value = 0.003589645367032649

$ws.Range("A9").Value = "Neutro"
$ws.Range("B9").Value = "Ltf"
$ws.Range("C9").Value = "Lrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.99918866666667
$ws.Range("H9").Value = 47.99756600000001
$ws.Range("I9").Value = 0.9848882775340866
$ws.Range("J9").Value = 0.9848882775340865
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 143.0355533333334
$ws.Range("N9").Value = 429.10666
$ws.Range("O9").Value = 0.2980616520156925
$ws.Range("P9").Value = 0.2980616520156925
$ws.Range("Q9").Value = 2288.452803821063
$ws.Range("R9").Value = 20596.07523438957
$ws.Range("S9").Value = 0.2935574270526997
$ws.Range("T9").Value = 0.2935574270526996

$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Ltf"
$ws.Range("C10").Value = "Lrp1"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.99918866666667
$ws.Range("H10").Value = 47.99756600000001
$ws.Range("I10").Value = 0.9848882775340866
$ws.Range("J10").Value = 0.9848882775340865
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 132.804812
$ws.Range("N10").Value = 398.414436
$ws.Range("O10").Value = 0.2767425352500014
$ws.Range("P10").Value = 0.2767425352500014
$ws.Range("Q10").Value = 2124.769243029198
$ws.Range("R10").Value = 19122.92318726278
$ws.Range("S10").Value = 0.2725604788627901
$ws.Range("T10").Value = 0.2725604788627901

$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Ltf"
$ws.Range("C11").Value = "Lrp1"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 15.99918866666667
$ws.Range("H11").Value = 47.99756600000001
$ws.Range("I11").Value = 0.9848882775340866
$ws.Range("J11").Value = 0.9848882775340865
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 125.707184
$ws.Range("N11").Value = 377.121552
$ws.Range("O11").Value = 0.261952291301752
$ws.Range("P11").Value = 0.261952291301752
$ws.Range("Q11").Value = 2011.212953571382
$ws.Range("R11").Value = 18100.91658214244
$ws.Range("S11").Value = 0.2579937409762899
$ws.Range("T11").Value = 0.2579937409762898

$ws.Range("A12").Value = "Neutro"
$ws.Range("B12").Value = "Ltf"
$ws.Range("C12").Value = "Lrp1"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 15.99918866666667
$ws.Range("H12").Value = 47.99756600000001
$ws.Range("I12").Value = 0.9848882775340866
$ws.Range("J12").Value = 0.9848882775340865
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 53.58648166666666
$ws.Range("N12").Value = 160.759445
$ws.Range("O12").Value = 0.111665071229204
$ws.Range("P12").Value = 0.1116650712292041
$ws.Range("Q12").Value = 857.3402301678746
$ws.Range("R12").Value = 7716.062071510871
$ws.Range("S12").Value = 0.1099776196636519
$ws.Range("T12").Value = 0.1099776196636519

$ws.Range("A13").Value = "Neutro"
$ws.Range("B13").Value = "Ltf"
$ws.Range("C13").Value = "Lrp1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 15.99918866666667
$ws.Range("H13").Value = 47.99756600000001
$ws.Range("I13").Value = 0.9848882775340866
$ws.Range("J13").Value = 0.9848882775340865
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 23.00271466666667
$ws.Range("N13").Value = 69.008144
$ws.Range("O13").Value = 0.04793372678759353
$ws.Range("P13").Value = 0.04793372678759355
$ws.Range("Q13").Value = 368.0247717975005
$ws.Range("R13").Value = 3312.222946177505
$ws.Range("S13").Value = 0.0472093656116225
$ws.Range("T13").Value = 0.04720936561162251
